$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Input" sheet: the three pending rows get re-populated with a new batch of
# data (CheA_mNG dataset), and the old 4th (already processed) row is removed.
# ---------------------------------------------------------------------------
$input = $wb.Worksheets.Item("Input")

$input.Cells.Item(1,1).Value = "1222 fliC- CheA_mNG"
$input.Cells.Item(1,2).Value = 20210519
$input.Cells.Item(1,2).NumberFormat = "0"
$input.Cells.Item(1,3).Value = "2s interval-2h37"

$input.Cells.Item(2,1).Value = "1222 fliC- CheA_mNG"
$input.Cells.Item(2,2).Value = 20210603
$input.Cells.Item(2,2).NumberFormat = "0"
$input.Cells.Item(2,3).Value = "2s interval-2h37"

$input.Cells.Item(3,1).Value = "1222 fliC- CheA_mNG"
$input.Cells.Item(3,2).Value = 20220825
$input.Cells.Item(3,3).Value = "2s interval-2h37"

# The former row 4 (20220825 entry) has now been filed away into the "CheA"
# sheet below, so remove it from the staging "Input" list.
$input.Rows(4).Delete()

$input.Range("A7").Activate()

# ---------------------------------------------------------------------------
# "CheA" sheet: file the newly-collected 2022-08-25 measurements (speed
# colourmap / polar loc vs plot runs) into the next free rows.
# ---------------------------------------------------------------------------
$cheA = $wb.Worksheets.Item("CheA")

$cheA.Cells.Item(7,1).Value = "1222 fliC- CheA_mNG"
$cheA.Cells.Item(7,2).Value = 20220825
$cheA.Cells.Item(7,2).NumberFormat = "0"
$cheA.Cells.Item(7,3).Value = "5s interval-1h37"

$cheA.Cells.Item(8,1).Value = "1222 fliC- CheA_mNG"
$cheA.Cells.Item(8,2).Value = 20220825
$cheA.Cells.Item(8,2).NumberFormat = "0"
$cheA.Cells.Item(8,3).Value = "2s interval-2h37"

$cheA.Cells.Item(9,1).Value = "1222 fliC- CheA_mNG"
$cheA.Cells.Item(9,2).Value = 20220825
$cheA.Cells.Item(9,2).NumberFormat = "0"
$cheA.Cells.Item(9,3).Value = "2s interval-3h37"

$cheA.Range("A3:C9").Select()
$cheA.Cells.Item(3,1).Activate()

$input.Activate()
